# feat: #108 Re-generate the statistics with the fixed minutes and seconds
# formatting in the haul fields.
#
# The "Общее время" (haul / total time) column (column I) contains strings
# like "17 ч. 3 мин. 24 сек." (H hours, M minutes, S seconds). This script
# zero-pads the minutes and seconds components to two digits, e.g.
# "17 ч. 3 мин. 24 сек." -> "17 ч. 03 мин. 24 сек." while leaving the
# hours component untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 9  # column I = "Общее время"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2].PadLeft(2, '0')
        $seconds = $matches[3].PadLeft(2, '0')
        $newVal = "$hours ч. $minutes мин. $seconds сек."

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
